# Updated tasks for Usecase: Show list of reportees
# Applies the edits described by the diff to "Klipper US.xlsx":
#   - Rewrites the task list on the "US 3" sheet (4th worksheet) with the new
#     scenario/task breakdown for the "Show list of reportees" use case.
#   - Removes the now-unused trailing blank rows (14-18) on that sheet.
#   - Widens column C on that sheet to fit the longer task descriptions.
#   - Moves the active/selected sheet & cell selections to match the new state.

$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Sheet "US List" (1st sheet): selection moves from G15 to D10, and it is no
# longer the active tab (handled automatically once we Activate "US 3" below).
# ------------------------------------------------------------------
$wsList = $wb.Worksheets.Item(1)
$wsList.Range("D10").Select()

# ------------------------------------------------------------------
# Sheet "US 3" (4th sheet): update the task breakdown table.
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(4)

# Row 5
$ws3.Range("C5").Value = "Identify roles"
$ws3.Range("D5").Value = 2
$ws3.Range("E5").Value = "Krutika"
$ws3.Range("F5").Value = "To do"

# Row 6
$ws3.Range("C6").Value = "Get list of reportees for individual Admin/Lead"
$ws3.Range("D6").ClearContents()
$ws3.Range("E6").Value = "Sidhdesh"
$ws3.Range("F6").Value = "Done"

# Row 7
$ws3.Range("C7").Value = "If role is Admin/Lead, then show list of reportees"
$ws3.Range("D7").Value = 3
$ws3.Range("E7").Value = "Krutika"
$ws3.Range("F7").Value = "To do"

# Row 8
$ws3.Range("C8").Value = "On selection of reportee, show one week Attendance record"
$ws3.Range("D8").Value = 3
$ws3.Range("E8").Value = "Krutika"
$ws3.Range("F8").Value = "To do"

# Row 9
$ws3.Range("C9").Value = "Check whether usecase is running with all scenarios"
$ws3.Range("D9").Value = 3
$ws3.Range("E9").Value = "Krutika"
$ws3.Range("F9").Value = "To do"

# Remove the now-empty trailing rows 14-18 (dimension shrinks to B2:F13).
$ws3.Range("B14:B18").EntireRow.Delete()

# Widen column C so the longer task text fits (stored width ends up at 56).
$ws3.Columns.Item(3).ColumnWidth = 55.16666666666666

# "US 3" becomes the active sheet/tab, with C11 selected.
$ws3.Activate()
$ws3.Range("C11").Select()
